# Move the "GitHub Link" slide (currently the 17th slide) so that it
# appears right after the "Conclusion" slide, i.e. becomes the 13th slide.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)
$s.MoveTo(13)
